$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("E2").Value = "SKIP"
$ws.Range("E11").Value = "SKIP"
$ws.Range("E12").Value = "SKIP"
